$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 212-213 (shifts old 212..243 down to 214..245),
# matching dimension change A1:R243 -> A1:R245.
$ws.Range("A212:A213").EntireRow.Insert()

# New row 212: Apio, Americana (o), Primera, date 2021-10-22 (44491)
$ws.Cells.Item(212, 1).Value = 8
$ws.Cells.Item(212, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(212, 3).Value = "Coquimbo"
$ws.Cells.Item(212, 4).Value = 44491
$ws.Cells.Item(212, 5).Value = 4
$ws.Cells.Item(212, 6).Value = 100112017
$ws.Cells.Item(212, 7).Value = "Apio"
$ws.Cells.Item(212, 8).Value = "Americana (o)"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 3300
$ws.Cells.Item(212, 11).Value = 6500
$ws.Cells.Item(212, 12).Value = 7000
$ws.Cells.Item(212, 13).Value = 6750
$ws.Cells.Item(212, 14).Value = "$/docena de matas"
$ws.Cells.Item(212, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(212, 16).Value = 1125
$ws.Cells.Item(212, 17).Value = 6
$ws.Cells.Item(212, 18).Value = "Hortaliza"

# New row 213: Apio, Americana (o), Segunda, date 2021-10-22 (44491)
$ws.Cells.Item(213, 1).Value = 8
$ws.Cells.Item(213, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = 44491
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = 100112017
$ws.Cells.Item(213, 7).Value = "Apio"
$ws.Cells.Item(213, 8).Value = "Americana (o)"
$ws.Cells.Item(213, 9).Value = "Segunda"
$ws.Cells.Item(213, 10).Value = 1200
$ws.Cells.Item(213, 11).Value = 5500
$ws.Cells.Item(213, 12).Value = 6000
$ws.Cells.Item(213, 13).Value = 5750
$ws.Cells.Item(213, 14).Value = "$/docena de matas"
$ws.Cells.Item(213, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(213, 16).Value = 958
$ws.Cells.Item(213, 17).Value = 6
$ws.Cells.Item(213, 18).Value = "Hortaliza"

# Make sure the Fecha column keeps the date style/format used elsewhere.
$ws.Range("D212:D213").NumberFormat = "YYYY-MM-DD HH:MM:SS"
